# The deck had a duplicated "Representação gráfica" slide (same image13.png
# picture placeholder repeated back-to-back). Remove the extra copy: this is
# slide id="360" which sits at position 12 in the deck (right after its
# identical twin at position 11, id="359").
$p = $ppt.ActivePresentation
$p.Slides.Item(12).Delete()
